$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the formatting (bold, border,
# centered) from the existing header cell H1, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-6
$dataI = @(1, 1, 9, 1, 3)
$dataJ = @(3, 6, 9, 3, 4)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
